# Adding logic for interventions of electrification
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("One Fed")

# Column E (Description of Measure) cells that move to the new, generic
# "Electrificaiton" description - this also causes the old, now-unused
# shared strings ("Electric Chillers ", "At bathroom upgrades",
# "AHU coils, and ASHP hot water for induction systems",
# "AHU coils, and elect HW boiler  for induction systems",
# "Electric boiler ") to be dropped from the workbook.
$ws.Range("E23").Value = "Electrificaiton"
$ws.Range("E24").Value = "Electrificaiton"
$ws.Range("E28").Value = "Electrificaiton"
$ws.Range("E29").Value = "Electrificaiton"
$ws.Range("E30").Value = "Electrificaiton"

# Updated "Change in Electricity Consumption Reduction (kWh)" (column I)
# and "Change in Steam Consumption, kLbs" (column K) figures for the
# electrification interventions.
$ws.Range("I23").Value = 3.5
$ws.Range("K23").Value = 3

$ws.Range("I28").Value = 2

$ws.Range("I29").Value = 1.05

$ws.Range("I30").Value = 1.05

# Restore the view to a clean scroll position and move the active
# selection to K28.
$ws.Activate()
$ws.Range("K28").Select()
